# Recreating the main files
# Applies the data refresh to the "morning"/"night" comparison table on the
# first worksheet: updated capacity/v-over-C/image-sum numbers, updated
# per-direction counts, several rows relabeled (B->C/A letter groups), and
# the now-unused tail rows (19-28 on the left, 18-27 on the right) cleared
# out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: capacity ---------------------------------------------------
$ws.Range("B2").Value = 1800
$ws.Range("F2").Value = 1800

# --- Row 3: v/C ----------------------------------------------------------
$ws.Range("B3").Value = 1.097241379310345
$ws.Range("F3").Value = 1.282068965517241

# --- Row 4: sum of images -------------------------------------------------
$ws.Range("B4").Value = 1591
$ws.Range("F4").Value = 1859

# --- Row 7: imageE ---------------------------------------------------------
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 0

# --- Row 8: imageD ---------------------------------------------------------
$ws.Range("B8").Value = 0
$ws.Range("F8").Value = 0

# --- Row 9: imageC ---------------------------------------------------------
$ws.Range("B9").Value = 594
$ws.Range("F9").Value = 615

# --- Row 10: imageB ----------------------------------------------------
$ws.Range("B10").Value = 76
$ws.Range("F10").Value = 110

# --- Row 11: imageA ----------------------------------------------------
$ws.Range("B11").Value = 921
$ws.Range("F11").Value = 1134

# --- Row 12 ----------------------------------------------------------------
$ws.Range("A12").Value = "CSt"
$ws.Range("B12").Value = 594
$ws.Range("E12").Value = "CNt"
$ws.Range("F12").Value = 615

# --- Row 13 ----------------------------------------------------------------
$ws.Range("A13").Value = "CNt"
$ws.Range("B13").Value = 594
$ws.Range("E13").Value = "CNl"
$ws.Range("F13").Value = 615

# --- Row 14 ----------------------------------------------------------------
$ws.Range("A14").Value = "BEr"
$ws.Range("B14").Value = 76
$ws.Range("E14").Value = "CEr"
$ws.Range("F14").Value = 312

# --- Row 15 ----------------------------------------------------------------
$ws.Range("A15").Value = "BEl"
$ws.Range("B15").Value = 76
$ws.Range("E15").Value = "BEl"
$ws.Range("F15").Value = 110

# --- Row 16 ----------------------------------------------------------------
$ws.Range("A16").Value = "ANt"
$ws.Range("B16").Value = 921
$ws.Range("E16").Value = "ASt"
$ws.Range("F16").Value = 774

# --- Row 17 ----------------------------------------------------------------
$ws.Range("A17").Value = "ANl"
$ws.Range("B17").Value = 560
$ws.Range("E17").Value = "ANt"
$ws.Range("F17").Value = 1134

# --- Row 18 ----------------------------------------------------------------
$ws.Range("A18").Value = "AEr"
$ws.Range("B18").Value = 609
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()

# --- Rows 19-27: left (A/B) and right (E/F) tails now empty ---------------
$ws.Range("A19:B27").ClearContents()
$ws.Range("E19:F27").ClearContents()

# --- Row 28: left tail (A/B) now empty; E/F were already empty -----------
$ws.Range("A28:B28").ClearContents()
